# Update "最後修改時間" (last modified time) values in column E
# to reflect the 2022/01/20 auto-commit refresh described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 5;   Value = "2022年01月20日 10:39:31" },
    @{ Row = 7;   Value = "2022年01月20日 10:39:54" },
    @{ Row = 9;   Value = "2022年01月20日 10:18:49" },
    @{ Row = 13;  Value = "2022年01月20日 10:19:22" },
    @{ Row = 19;  Value = "2022年01月20日 10:22:23" },
    @{ Row = 22;  Value = "2022年01月20日 10:23:13" },
    @{ Row = 25;  Value = "2022年01月20日 10:24:19" },
    @{ Row = 26;  Value = "2022年01月20日 10:24:44" },
    @{ Row = 28;  Value = "2022年01月20日 10:25:35" },
    @{ Row = 39;  Value = "2022年01月20日 10:36:21" },
    @{ Row = 41;  Value = "2022年01月20日 10:37:15" },
    @{ Row = 44;  Value = "2022年01月20日 10:37:42" },
    @{ Row = 49;  Value = "2022年01月20日 10:38:18" },
    @{ Row = 53;  Value = "2022年01月20日 10:38:57" },
    @{ Row = 54;  Value = "2022年01月20日 10:39:16" },
    @{ Row = 55;  Value = "2022年01月20日 10:42:24" },
    @{ Row = 56;  Value = "2022年01月20日 10:43:18" },
    @{ Row = 57;  Value = "2022年01月20日 10:44:03" },
    @{ Row = 58;  Value = "2022年01月20日 10:44:20" },
    @{ Row = 59;  Value = "2022年01月20日 10:44:37" },
    @{ Row = 60;  Value = "2022年01月20日 10:44:53" },
    @{ Row = 61;  Value = "2022年01月20日 10:45:09" },
    @{ Row = 64;  Value = "2022年01月20日 10:48:13" },
    @{ Row = 65;  Value = "2022年01月20日 10:49:14" },
    @{ Row = 66;  Value = "2022年01月20日 10:57:04" },
    @{ Row = 67;  Value = "2022年01月20日 10:55:17" },
    @{ Row = 68;  Value = "2022年01月20日 10:58:09" },
    @{ Row = 69;  Value = "2022年01月20日 10:59:19" },
    @{ Row = 70;  Value = "2022年01月20日 10:59:49" },
    @{ Row = 71;  Value = "2022年01月20日 11:00:17" },
    @{ Row = 73;  Value = "2022年01月20日 11:00:42" },
    @{ Row = 74;  Value = "2022年01月20日 11:00:58" },
    @{ Row = 75;  Value = "2022年01月20日 11:01:44" },
    @{ Row = 76;  Value = "2022年01月20日 11:02:26" },
    @{ Row = 77;  Value = "2022年01月20日 11:03:12" },
    @{ Row = 80;  Value = "2022年01月20日 11:03:31" },
    @{ Row = 81;  Value = "2022年01月20日 11:03:54" },
    @{ Row = 83;  Value = "2022年01月20日 11:05:21" },
    @{ Row = 84;  Value = "2022年01月20日 11:06:26" },
    @{ Row = 85;  Value = "2022年01月20日 11:07:17" },
    @{ Row = 86;  Value = "2022年01月20日 11:07:45" },
    @{ Row = 87;  Value = "2022年01月20日 11:19:19" },
    @{ Row = 89;  Value = "2022年01月20日 11:20:01" },
    @{ Row = 92;  Value = "2022年01月20日 11:20:23" },
    @{ Row = 100; Value = "2022年01月20日 11:21:28" },
    @{ Row = 103; Value = "2022年01月20日 11:22:18" },
    @{ Row = 112; Value = "2022年01月20日 11:24:01" },
    @{ Row = 113; Value = "2022年01月20日 11:26:29" },
    @{ Row = 115; Value = "2022年01月20日 11:26:55" },
    @{ Row = 118; Value = "2022年01月20日 11:27:22" },
    @{ Row = 121; Value = "2022年01月20日 11:28:00" },
    @{ Row = 122; Value = "2022年01月20日 11:28:25" },
    @{ Row = 123; Value = "2022年01月20日 11:28:43" },
    @{ Row = 124; Value = "2022年01月20日 11:29:25" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.Value
}
